# Update the twenty-five division problems scattered through the worksheet
# table. Every "before" value in the table is unique, so each one can be
# located and swapped with Find/Replace ("wrap find" off, "replace: all" = 2)
# without touching the blank practice rows in between.
#
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)
#
# NOTE: "71÷2=" is both an old value (row 3, col 5) and a new value
# (row 4, col 5 replacing "19÷3="). The "40÷4=" ... "71÷2="->"47÷6="
# replacement below MUST run before the "19÷3="->"71÷2=" replacement,
# otherwise the later Find would match two cells instead of one. The
# document order used here already satisfies that constraint.
$d = $word.ActiveDocument

$d.Content.Find.Execute("40÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷4=", 2) | Out-Null
$d.Content.Find.Execute("36÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷6=", 2) | Out-Null
$d.Content.Find.Execute("52÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷5=", 2) | Out-Null
$d.Content.Find.Execute("96÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷2=", 2) | Out-Null
$d.Content.Find.Execute("37÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷2=", 2) | Out-Null
$d.Content.Find.Execute("21÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷7=", 2) | Out-Null
$d.Content.Find.Execute("61÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷5=", 2) | Out-Null
$d.Content.Find.Execute("80÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷9=", 2) | Out-Null
$d.Content.Find.Execute("94÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "50÷6=", 2) | Out-Null
$d.Content.Find.Execute("69÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷4=", 2) | Out-Null
$d.Content.Find.Execute("97÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷8=", 2) | Out-Null
$d.Content.Find.Execute("18÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷7=", 2) | Out-Null
$d.Content.Find.Execute("54÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷3=", 2) | Out-Null
$d.Content.Find.Execute("25÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷6=", 2) | Out-Null
$d.Content.Find.Execute("71÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷6=", 2) | Out-Null
$d.Content.Find.Execute("91÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=", 2) | Out-Null
$d.Content.Find.Execute("53÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷4=", 2) | Out-Null
$d.Content.Find.Execute("65÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷9=", 2) | Out-Null
$d.Content.Find.Execute("37÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "38÷7=", 2) | Out-Null
$d.Content.Find.Execute("19÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=", 2) | Out-Null
$d.Content.Find.Execute("51÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷3=", 2) | Out-Null
$d.Content.Find.Execute("27÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷5=", 2) | Out-Null
$d.Content.Find.Execute("49÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷8=", 2) | Out-Null
$d.Content.Find.Execute("10÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=", 2) | Out-Null
$d.Content.Find.Execute("90÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=", 2) | Out-Null
